$wb = $excel.ActiveWorkbook
$wsJob = $wb.Worksheets.Item(1)   # "Job Type"
$wsTime = $wb.Worksheets.Item(2)  # "Processing Time"

# ---------------------------------------------------------------------------
# Sheet "Job Type": update mutated Job Type / setup values for J4-J10, fix J11
# and append new rows for J12-J20
# ---------------------------------------------------------------------------
$wsJob.Range("B5").Value = 1
$wsJob.Range("C5").Value = 6
$wsJob.Range("B6").Value = 1
$wsJob.Range("C6").Value = 6

$wsJob.Range("B9").Value = 2
$wsJob.Range("C9").Value = 8
$wsJob.Range("B10").Value = 2
$wsJob.Range("C10").Value = 8
$wsJob.Range("B11").Value = 2
$wsJob.Range("C11").Value = 8

$wsJob.Range("C12").Value = 3

$jobTypeNewRows = @(
    @{ Row = 13; Name = "J12"; B = 3; C = 3 },
    @{ Row = 14; Name = "J13"; B = 3; C = 3 },
    @{ Row = 15; Name = "J14"; B = 3; C = 3 },
    @{ Row = 16; Name = "J15"; B = 3; C = 3 },
    @{ Row = 17; Name = "J16"; B = 4; C = 4 },
    @{ Row = 18; Name = "J17"; B = 4; C = 4 },
    @{ Row = 19; Name = "J18"; B = 4; C = 4 },
    @{ Row = 20; Name = "J19"; B = 4; C = 4 },
    @{ Row = 21; Name = "J20"; B = 4; C = 4 }
)

foreach ($r in $jobTypeNewRows) {
    $wsJob.Cells.Item($r.Row, 1).Value = $r.Name
    $wsJob.Cells.Item($r.Row, 1).HorizontalAlignment = -4108
    $wsJob.Cells.Item($r.Row, 1).VerticalAlignment = -4108
    $wsJob.Cells.Item($r.Row, 2).Value = $r.B
    $wsJob.Cells.Item($r.Row, 2).HorizontalAlignment = -4108
    $wsJob.Cells.Item($r.Row, 2).VerticalAlignment = -4108
    $wsJob.Cells.Item($r.Row, 3).Value = $r.C
    $wsJob.Cells.Item($r.Row, 3).HorizontalAlignment = -4108
    $wsJob.Cells.Item($r.Row, 3).VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# Sheet "Processing Time": fix processing-time bug for J3-J11 and append
# new rows for J12-J20
# ---------------------------------------------------------------------------
$wsTime.Range("B4").Value = 12
$wsTime.Range("B5").Value = 12
$wsTime.Range("B6").Value = 15
$wsTime.Range("B7").Value = 4
$wsTime.Range("B8").Value = 6
$wsTime.Range("B9").Value = 7
$wsTime.Range("B10").Value = 9
$wsTime.Range("B11").Value = 10
$wsTime.Range("B12").Value = 2

$processingNewRows = @(
    @{ Row = 13; Name = "J12"; B = 3 },
    @{ Row = 14; Name = "J13"; B = 3 },
    @{ Row = 15; Name = "J14"; B = 4 },
    @{ Row = 16; Name = "J15"; B = 5 },
    @{ Row = 17; Name = "J16"; B = 8 },
    @{ Row = 18; Name = "J17"; B = 12 },
    @{ Row = 19; Name = "J18"; B = 13 },
    @{ Row = 20; Name = "J19"; B = 15 },
    @{ Row = 21; Name = "J20"; B = 17 }
)

foreach ($r in $processingNewRows) {
    $wsTime.Cells.Item($r.Row, 1).Value = $r.Name
    $wsTime.Cells.Item($r.Row, 1).HorizontalAlignment = -4108
    $wsTime.Cells.Item($r.Row, 1).VerticalAlignment = -4108
    $wsTime.Cells.Item($r.Row, 2).Value = $r.B
    $wsTime.Cells.Item($r.Row, 2).HorizontalAlignment = -4108
    $wsTime.Cells.Item($r.Row, 2).VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# Selections / active sheet: select a range on "Processing Time" first, then
# finish on "Job Type" so that "Job Type" ends up as the active/visible tab
# (matches tabSelected moving from sheet2 to sheet1).
# ---------------------------------------------------------------------------
$wsTime.Range("A13:A21").Select()
$wsJob.Range("D21").Select()
